$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update week 5 attendance values (rows 147-151, columns D/E/G/H = Walter/Paige/Benjamin/Hadewij)
$ws.Range("D147").Value = 4
$ws.Range("E147").Value = 4
$ws.Range("G147").Value = 4
$ws.Range("H147").Value = 4

$ws.Range("D148").Value = 6
$ws.Range("G148").Value = 6
$ws.Range("H148").Value = 6

$ws.Range("D149").Value = 2
$ws.Range("E149").Value = 2
$ws.Range("G149").Value = 2
$ws.Range("H149").Value = 2

$ws.Range("D150").Value = 4
$ws.Range("E150").Value = 4
$ws.Range("G150").Value = 4
$ws.Range("H150").Value = 4

$ws.Range("B151").Value = 0
$ws.Range("I151").Value = 0

$ws.Range("B152").Value = 16

# Update selection / view to match the saved state
$ws.Range("E148").Select()
